# Add new columns I ("I0") and J ("IF") to the worksheet, matching the
# existing header style (bold, centered, thin-bordered) used by column H,
# and fill in the per-row numeric values for rows 2-63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the formatting of the existing "IP" header cell (H1) onto the two
# new header cells so they pick up the same style index (bold font,
# centered/top alignment, thin border) as every other header cell.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(1, 9).Value  = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# --- Data ----------------------------------------------------------------
# Values for column I (I0) and column J (IF) for rows 2 through 63,
# in row order.
$iValues = @(3,6,8,5,7,7,8,7,7,6,6,9,5,9,7,6,7,9,5,7,4,7,9,7,9,7,9,5,8,6,8,4,9,7,6,7,8,8,9,5,8,8,8,7,4,6,10,6,6,6,5,6,7,8,6,6,5,5,9,4,9,7)
$jValues = @(4,6,9,6,7,7,8,7,7,6,6,9,5,9,7,7,7,9,6,7,5,8,9,7,9,7,9,6,8,6,8,5,9,7,6,7,8,8,9,6,8,8,8,8,6,7,11,7,6,7,7,7,8,8,6,7,6,5,9,4,9,7)

for ($r = 2; $r -le 63; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value  = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}

Write-Host "Added columns I (I0) and J (IF) with data for rows 1-63"
